# chore(auto): pipeline update & logic validation
#
# The war-participation data pipeline had produced rows whose "Player
# Status" / war-score columns (B, D:H) were offset from the correct
# player row. This re-applies the corrected (shifted) values for the
# affected rows, while leaving the player names (column A) and the
# "Fonte de Dados" column (C) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-WarRow($Row, $Status, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 2).Value = $Status   # B: Player Status
    $ws.Cells.Item($Row, 4).Value = $D         # D: Ultima Guerra
    $ws.Cells.Item($Row, 5).Value = $E         # E: Guerra -2
    $ws.Cells.Item($Row, 6).Value = $F         # F: Guerra -3
    $ws.Cells.Item($Row, 7).Value = $G         # G: Guerra -4
    $ws.Cells.Item($Row, 8).Value = $H         # H: Guerra -5
}

# Rows 24-25 (StelaAby / ^_^^_^) swap their war-data between each other
Set-WarRow 24 "Razoável"  12 16 12 12 16
Set-WarRow 25 "Verificar" 0  0  0  0  0

# Rows 47-58 cycle their war-data up by one row (47<-48<-49<-...<-58<-47)
Set-WarRow 47 "Verificar" 0  0  0  0  0
Set-WarRow 48 "Ok"        16 16 12 16 13
Set-WarRow 49 "Ok"        16 16 0  0  16
Set-WarRow 50 "Razoável"  15 16 16 16 16
Set-WarRow 51 "Ok"        16 16 16 16 16
Set-WarRow 52 "Ok"        16 16 16 16 16
Set-WarRow 53 "Ok"        16 16 12 15 16
Set-WarRow 54 "Ok"        16 16 16 16 16
Set-WarRow 55 "Ok"        16 16 16 16 16
Set-WarRow 56 "Verificar" 8  14 7  10 6
Set-WarRow 57 "Ok"        16 16 16 16 16
Set-WarRow 58 "Ok"        16 14 16 16 15
